# Oregon_Converted.xlsx update — "Updated policies and graphs"
#
# 1. Row 7 ("Weights"): the Q (Customer Mandate) weight is zeroed out and the
#    total weight (Z7) drops from 13 to 12.
# 2. Rows 20-221: the per-day "LockdownEffectiveness" (Z) column is
#    recomputed against the new weight total.
# 3. Twelve new days (9/30/2020 - 10/11/2020) are appended as rows 222-233.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- Row 7: Weights row changes (Q7 and Z7) ---
$ws.Cells.Item(7,17).Value = 0    # Q7: 1 -> 0
$ws.Cells.Item(7,26).Value = 12   # Z7: 13 -> 12

# --- Rows 20-221: recomputed LockdownEffectiveness (Z) values ---
for ($r = 20; $r -le 23; $r++) { $ws.Cells.Item($r,26).Value = 0.08333333333333333 }
$ws.Cells.Item(24,26).Value = 0.1666666666666667
for ($r = 25; $r -le 30; $r++) { $ws.Cells.Item($r,26).Value = 0.25 }
for ($r = 31; $r -le 74; $r++) { $ws.Cells.Item($r,26).Value = 0.8333333333166668 }
for ($r = 75; $r -le 84; $r++) { $ws.Cells.Item($r,26).Value = 0.7777777777666669 }
for ($r = 85; $r -le 100; $r++) { $ws.Cells.Item($r,26).Value = 0.5972222222166667 }
for ($r = 101; $r -le 108; $r++) { $ws.Cells.Item($r,26).Value = 0.6805555555499999 }
for ($r = 109; $r -le 119; $r++) { $ws.Cells.Item($r,26).Value = 0.5694444444416666 }
for ($r = 120; $r -le 221; $r++) { $ws.Cells.Item($r,26).Value = 0.4166666666666667 }

# --- New rows 222-233: additional days 9/30/2020 - 10/11/2020 ---
$newDates = @("9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020")

$rowValues = @(0,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,0,0,0,0)

$destRow = 222
foreach ($d in $newDates) {
    # Clone A221's style (bold/border/center, General fmt, shared-string type)
    # into the new row's A cell so formatting matches the rest of the date column.
    $ws.Range("A221").Copy($ws.Cells.Item($destRow, 1))

    # Write the literal date text through a helper formula cell + paste-values,
    # so Excel doesn't auto-convert the date-look-alike string into a real date.
    $ws.Range("AA1").Formula = '="' + $d + '"'
    $ws.Range("AA1").Copy()
    $ws.Cells.Item($destRow, 1).PasteSpecial(-4163)
    $ws.Range("AA1").Clear()

    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($destRow, $col + 2).Value = $rowValues[$col]
    }
    $ws.Cells.Item($destRow, 26).Value = 0.4166666666666667

    $destRow = $destRow + 1
}
